# Apply the "Upload new version with timestamp" edit:
#  - CETAL 500MG 20 TAB row (row 11): unit price doubled (12.0000 -> 24.0000)
#    and the two ratio/traffic columns swapped (H11 <-> Q11 values)
#  - Grand total (P28) increases to match the new price (809.325 -> 821.325)
#  - Footer timestamp (A29) updated to the new export time (12:49 PM -> 12:52 PM)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H11").Value = "0:1"
$ws.Range("Q11").Value = "1:0"

# P11 is numeric-formatted ("0.00") but must keep storing a literal text
# value (as in the source file). Temporarily flip the cell to a text
# format so Excel doesn't coerce the assignment into a number, then
# restore the original number format so the cell style is unchanged.
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "24.0000"
$ws.Range("P11").NumberFormat = "0.00"

$ws.Range("P28").Value = 821.32500000000005

$ws.Range("A29").Value = "Tuesday, 17 June, 2025 12:52 PM"
